$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 13.008
$ws.Range("C10").Value = -12.214
$ws.Range("C12").Value = -12.53
$ws.Range("D13").Value = -7.831999999999999
$ws.Range("C18").Value = -12.283
$ws.Range("E20").Value = 12.932
